# Adding CCES NBHF detections to Tethys
# NBHF detections are saved as "Odontoceti, Group: NBHG"
# Updated detection worksheets to upload to Tethys

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: drop "Parameter 2".."Parameter 6" (now-unused extra parameter
# columns), keep the bold header style on G1:K1 but clear their text.
$ws.Range("G1:K1").ClearContents()

# Row 2: rename the NBHF detection parameters.
#   F2 "min" -> "nClicks"        (write this first so it lands in the
#                                 shared-string table before "NBHF")
#   B2 "NBHF Odontocete" -> "NBHF"
#   C2 "UO" -> "NBHF"
$ws.Range("F2").Value = "nClicks"
$ws.Range("B2").Value = "NBHF"
$ws.Range("C2").Value = "NBHF"

# Row 2: drop the now-unused "max"/"peak 1..4" parameter cells.
$ws.Range("G2:K2").ClearContents()

# Column C widened to fit the new "NBHF" values (best-fit autosize).
$ws.Columns("C:C").ColumnWidth = 11.9453125

# Selection moves up from B3 to B2.
$ws.Range("B2").Select()
